# Remove the blank spacer row (row 3) so the second author's data moves up to row 3.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(3).Delete()

# Underline the hyperlink font used by the e-mail cell style (style s="2").
$ws.Range("E2").Font.Underline = $true

# Underline the header-ish row (row 2) fields, matching the new style s="1" used by the template.
$row2Ranges = @("A2","B2","C2","D2","F2","G2","H2","I2","J2","K2","L2","M2","N2","Q2","R2")
foreach ($addr in $row2Ranges) {
    $ws.Range($addr).Font.Underline = $true
}

# The (new) row 3, first cell also picks up the underline formatting.
$ws.Range("A3").Font.Underline = $true

# Update the active selection to match the edited sheet.
$null = $ws.Range("A3").Select()
